$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.734.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.441.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.03%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.26"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.85%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.435.34"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.53%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.13"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.02%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.989.42"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.84%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.445.51"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.09%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.733.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.54%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.71%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.45%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +19.01%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.31%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.81"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.11%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.14"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.75%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "64.48"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.83%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "594.03"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.112"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.06%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.31%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.74%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.386"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.85%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.31"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.67%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0771"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.187.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0430"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.38%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.81"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +22.88%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.27%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.35%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.68%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.72"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.85%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.11"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.02%  "

